$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 33) following the existing table pattern.
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Update the active selection / view as captured in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("E29").Select()
